# Phase 1 now emits an equipment list - add a "Modem Name" column to the
# Outstation worklist template, between "Oustation Specific Model" (R) and
# "Modem Manufacturer" (S), populating every data row with "Modem".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at S; this shifts the existing Modem*/Comment columns
# one position to the right (S->T, T->U, ... Y->Z) and copies the left
# neighbour's (R) formatting, matching Excel's native "Insert" behaviour.
$ws.Columns("S:S").Insert()

# The new column inherits column R's width - set it to match.
$ws.Range("S1").ColumnWidth = $ws.Range("R1").ColumnWidth

# Header for the new column.
$ws.Range("S1").Value2 = "Modem Name"

# Populate every data row (2-11) with the new "Modem" value.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 19).Value2 = "Modem"
}

# Reflect the user's scrolling/selection while editing near column J.
$null = $ws.Range("J11").Select()
